$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 29-30, pushing existing data (rows 29..113) down to 31..115
$ws.Rows("29:30").Insert()

# Fill in row 29 (new record)
$ws.Range("A29").Value = 11
$ws.Range("B29").Value = "Vega Monumental Concepción"
$ws.Range("C29").Value = "Bíobío"
$ws.Range("D29").Value = 44987
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100103
$ws.Range("H29").Value = "Frutos de hueso (carozo)"
$ws.Range("I29").Value = 100103002
$ws.Range("J29").Value = "Ciruela"
$ws.Range("K29").Value = "Black Amber"
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = 12000
$ws.Range("O29").Value = 12000
$ws.Range("P29").Value = 12000
$ws.Range("Q29").Value = '$/bandeja 18 kilos granel'
$ws.Range("R29").Value = "Región del Maule"
$ws.Range("S29").Value = 667
$ws.Range("T29").Value = 18

# Fill in row 30 (new record)
$ws.Range("A30").Value = 11
$ws.Range("B30").Value = "Vega Monumental Concepción"
$ws.Range("C30").Value = "Bíobío"
$ws.Range("D30").Value = 44987
$ws.Range("E30").Value = 8
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100103
$ws.Range("H30").Value = "Frutos de hueso (carozo)"
$ws.Range("I30").Value = 100103002
$ws.Range("J30").Value = "Ciruela"
$ws.Range("K30").Value = "Black Amber"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 10000
$ws.Range("O30").Value = 10000
$ws.Range("P30").Value = 10000
$ws.Range("Q30").Value = '$/bandeja 18 kilos granel'
$ws.Range("R30").Value = "Región del Maule"
$ws.Range("S30").Value = 556
$ws.Range("T30").Value = 18

# Apply same date style (style index 2 = YYYY-MM-DD HH:MM:SS) used by column D elsewhere
$ws.Range("D29").NumberFormat = $ws.Range("D31").NumberFormat
$ws.Range("D30").NumberFormat = $ws.Range("D31").NumberFormat
